# #249: added test case and updated test file
#
# Adds a small 4-column/10-row demo table ("Table1") to Sheet1, starting at
# O4, with header labels Column1..Column4 — mirroring the xlsx fixture used
# by the NvPr test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new table (O4:R4) -> becomes shared-string entries
# "Column1".."Column4".
$ws.Range("O4").Value = "Column1"
$ws.Range("P4").Value = "Column2"
$ws.Range("Q4").Value = "Column3"
$ws.Range("R4").Value = "Column4"

# Match the fixture's explicit column widths for the new columns (O:R).
$ws.Columns("O:R").ColumnWidth = 10.17

# Turn O4:R13 into a native Excel table (xlSrcRange=1, xlYes=1 for headers).
$tbl = $ws.ListObjects.Add(1, $ws.Range("O4:R13"), $null, 1)
$tbl.Name = "Table1"
$tbl.Comment = "Luke, I am your father... seriously..."

# Leave the selection where the author left it in the fixture.
$ws.Range("R24").Select() | Out-Null
